$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.75242533404057
$ws.Range("D2").Value = 7.998920270634828
$ws.Range("E2").Value = 10.32980958153274
$ws.Range("F2").Value = 41.96376400327019
$ws.Range("G2").Value = 3.709720673870575
$ws.Range("I2").Value = 33.08355626267235
$ws.Range("K2").Value = 13.29379090718554
$ws.Range("L2").Value = 10.56304082407664
$ws.Range("B3").Value = 12.72039652785334
$ws.Range("D3").Value = 7.992164242515012
$ws.Range("E3").Value = 10.28483599605408
$ws.Range("F3").Value = 41.51039658091621
$ws.Range("G3").Value = 3.713441198970362
$ws.Range("I3").Value = 32.95916790202578
$ws.Range("K3").Value = 13.08001272797512
$ws.Range("L3").Value = 10.52933145534091
$ws.Range("B4").Value = 12.70442667742202
$ws.Range("D4").Value = 7.988212781880506
$ws.Range("E4").Value = 10.2564717128333
$ws.Range("F4").Value = 41.23773778493378
$ws.Range("G4").Value = 3.715843385999123
$ws.Range("I4").Value = 32.88609577997138
$ws.Range("K4").Value = 12.95202536383974
$ws.Range("L4").Value = 10.51106533763075
$ws.Range("B5").Value = 12.69885524639157
$ws.Range("D5").Value = 7.986651754929225
$ws.Range("E5").Value = 10.24472396462557
$ws.Range("F5").Value = 41.12816116046846
$ws.Range("G5").Value = 3.716852023794093
$ws.Range("I5").Value = 32.85715875409383
$ws.Range("K5").Value = 12.90076995908587
$ws.Range("L5").Value = 10.50423802787787
$ws.Range("B6").Value = 12.69798685649458
$ws.Range("D6").Value = 7.986395512947109
$ws.Range("E6").Value = 10.24276182875473
$ws.Range("F6").Value = 41.11006149484344
$ws.Range("G6").Value = 3.71702130609005
$ws.Range("I6").Value = 32.85240481632607
$ws.Range("K6").Value = 12.89231565255234
$ws.Range("L6").Value = 10.50314170805159
$ws.Range("B7").Value = 12.70434773916439
$ws.Range("D7").Value = 7.988191530373823
$ws.Range("E7").Value = 10.25631404459326
$ws.Range("F7").Value = 41.23625365820091
$ws.Range("G7").Value = 3.715856868329509
$ws.Range("I7").Value = 32.88570210983634
$ws.Range("K7").Value = 12.95133036856567
$ws.Range("L7").Value = 10.51097076098035
$ws.Range("B8").Value = 12.74061923750409
$ws.Range("D8").Value = 7.996549802875549
$ws.Range("E8").Value = 10.31445817153619
$ws.Range("F8").Value = 41.80631751387369
$ws.Range("G8").Value = 3.710979135990818
$ws.Range("I8").Value = 33.03998787109376
$ws.Range("K8").Value = 13.21944771777121
$ws.Range("L8").Value = 10.55091627289668
$ws.Range("B9").Value = 12.84071634718517
$ws.Range("D9").Value = 8.014513289999442
$ws.Range("E9").Value = 10.42254628392938
$ws.Range("F9").Value = 42.9650698425578
$ws.Range("G9").Value = 3.702343211932766
$ws.Range("I9").Value = 33.36827869891567
$ws.Range("K9").Value = 13.76773489253133
$ws.Range("L9").Value = 10.6482994083494
$ws.Range("B10").Value = 12.93137159490391
$ws.Range("D10").Value = 8.028686474330735
$ws.Range("E10").Value = 10.49837599569955
$ws.Range("F10").Value = 43.83519897403432
$ws.Range("G10").Value = 3.696557704532085
$ws.Range("I10").Value = 33.624488610851
$ws.Range("K10").Value = 14.17962518738002
$ws.Range("L10").Value = 10.73109137072144
$ws.Range("B11").Value = 12.9761923642634
$ws.Range("D11").Value = 8.035349304271994
$ws.Range("E11").Value = 10.53210246473562
$ws.Range("F11").Value = 44.23379743398274
$ws.Range("G11").Value = 3.694045633341382
$ws.Range("I11").Value = 33.74414436287958
$ws.Range("K11").Value = 14.3679995448055
$ws.Range("L11").Value = 10.77110529689645
$ws.Range("B12").Value = 12.99366688039601
$ws.Range("D12").Value = 8.037903519947568
$ws.Range("E12").Value = 10.54476371879815
$ws.Range("F12").Value = 44.38502313673447
$ws.Range("G12").Value = 3.693111483840089
$ws.Range("I12").Value = 33.78988616922742
$ws.Range("K12").Value = 14.43939759657026
$ws.Range("L12").Value = 10.78658694701312
$ws.Range("B13").Value = 12.98988132789044
$ws.Range("D13").Value = 8.037352036035026
$ws.Range("E13").Value = 10.54204178831806
$ws.Range("F13").Value = 44.35244309474866
$ws.Range("G13").Value = 3.693311910059665
$ws.Range("I13").Value = 33.78001595144934
$ws.Range("K13").Value = 14.42401910585182
$ws.Range("L13").Value = 10.7832382008324
$ws.Range("B14").Value = 12.9776200232325
$ws.Range("D14").Value = 8.035558815965947
$ws.Range("E14").Value = 10.53314632094194
$ws.Range("F14").Value = 44.24623379706535
$ws.Range("G14").Value = 3.693968437857349
$ws.Range("I14").Value = 33.74789902790995
$ws.Range("K14").Value = 14.37387265349099
$ws.Range("L14").Value = 10.77237243590356
$ws.Range("B15").Value = 12.97017457393681
$ws.Range("D15").Value = 8.034464478011953
$ws.Range("E15").Value = 10.52768324531076
$ws.Range("F15").Value = 44.18121130878484
$ws.Range("G15").Value = 3.694372806086342
$ws.Range("I15").Value = 33.72828210624082
$ws.Range("K15").Value = 14.34316262319826
$ws.Range("L15").Value = 10.76575944519564
$ws.Range("B16").Value = 12.92851341432033
$ws.Range("D16").Value = 8.028255392072847
$ws.Range("E16").Value = 10.4961564487399
$ws.Range("F16").Value = 43.80919569335384
$ws.Range("G16").Value = 3.696724275465949
$ws.Range("I16").Value = 33.61672997619719
$ws.Range("K16").Value = 14.16732781568049
$ws.Range("L16").Value = 10.72852293648538
$ws.Range("B17").Value = 12.90386382584117
$ws.Range("D17").Value = 8.024501612709299
$ws.Range("E17").Value = 10.47661909514002
$ws.Range("F17").Value = 43.58160526305215
$ws.Range("G17").Value = 3.698197428902162
$ws.Range("I17").Value = 33.54908060968128
$ws.Range("K17").Value = 14.05965848637533
$ws.Range("L17").Value = 10.70627566046727
$ws.Range("B18").Value = 12.89002395912529
$ws.Range("D18").Value = 8.022362801053804
$ws.Range("E18").Value = 10.46530933860346
$ws.Range("F18").Value = 43.45096899050867
$ws.Range("G18").Value = 3.699056028950161
$ws.Range("I18").Value = 33.51046368781064
$ws.Range("K18").Value = 13.99782899488349
$ws.Range("L18").Value = 10.69370151418284
$ws.Range("B19").Value = 12.88539644248842
$ws.Range("D19").Value = 8.021642110717639
$ws.Range("E19").Value = 10.461467606727
$ws.Range("F19").Value = 43.40678724443384
$ws.Range("G19").Value = 3.69934867700652
$ws.Range("I19").Value = 33.4974394846152
$ws.Range("K19").Value = 13.97691406889166
$ws.Range("L19").Value = 10.68948249721756
$ws.Range("B20").Value = 12.90645293501308
$ws.Range("D20").Value = 8.024899110466087
$ws.Range("E20").Value = 10.47870636921605
$ws.Range("F20").Value = 43.60580576916279
$ws.Range("G20").Value = 3.698039442428382
$ws.Range("I20").Value = 33.5562517657818
$ws.Range("K20").Value = 14.07111040222349
$ws.Range("L20").Value = 10.70862101796079
$ws.Range("B21").Value = 12.98120795387306
$ws.Range("D21").Value = 8.036084681545287
$ws.Range("E21").Value = 10.53576212423808
$ws.Range("F21").Value = 44.27742320545989
$ws.Range("G21").Value = 3.693775136024263
$ws.Range("I21").Value = 33.75732098441163
$ws.Range("K21").Value = 14.3886007449029
$ws.Range("L21").Value = 10.77555511243842
$ws.Range("B22").Value = 13.03298326323099
$ws.Range("D22").Value = 8.043576714168344
$ws.Range("E22").Value = 10.57240902292871
$ws.Range("F22").Value = 44.7179739074428
$ws.Range("G22").Value = 3.691087888589482
$ws.Range("I22").Value = 33.89123442159154
$ws.Range("K22").Value = 14.59644164520602
$ws.Range("L22").Value = 10.82121499779597
$ws.Range("B23").Value = 13.00508730187217
$ws.Range("D23").Value = 8.039561404016538
$ws.Range("E23").Value = 10.5529085467814
$ws.Range("F23").Value = 44.48273368643041
$ws.Range("G23").Value = 3.69251303387552
$ws.Range("I23").Value = 33.81953863168665
$ws.Range("K23").Value = 14.48550726385592
$ws.Range("L23").Value = 10.7966733507533
$ws.Range("B24").Value = 12.90528136519251
$ws.Range("D24").Value = 8.024719341671879
$ws.Range("E24").Value = 10.47776295393831
$ws.Range("F24").Value = 43.59486406230839
$ws.Range("G24").Value = 3.69811083177269
$ws.Range("I24").Value = 33.55300882649713
$ws.Range("K24").Value = 14.06593276373304
$ws.Range("L24").Value = 10.7075600080392
$ws.Range("B25").Value = 12.81058856150726
$ws.Range("D25").Value = 8.009483524422983
$ws.Range("E25").Value = 10.39393716889066
$ws.Range("F25").Value = 42.64784728152466
$ws.Range("G25").Value = 3.704580718540229
$ws.Range("I25").Value = 33.27677821851757
$ws.Range("K25").Value = 13.61748162892115
$ws.Range("L25").Value = 10.61995071479073
